$wb = $excel.ActiveWorkbook

# The workbook's second tab ("FirstSet") is the active sheet and holds the
# credentials table whose password value needs updating.
$ws = $wb.Worksheets.Item("FirstSet")
$ws.Activate()

# Update the password in B2 from "Password2!" to "Password0!"
$ws.Range("B2").Value = "Password0!"

# Move/record the active selection to B2 (matches the saved sheetView state)
$ws.Range("B2").Select()
